# Remove the Lassa L-protein "Biochemical characterization" record (RefID 50)
# which was picked up again via a linked/duplicate PubMed search, and keep the
# rest of the reference list contiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the whole row first (mirrors the manual "right click row header -> Delete"
# workflow) so the sheet's remembered selection matches the target state.
$ws.Rows.Item(19).Select()
$ws.Rows.Item(19).Delete()

# Re-apply the AutoFilter over the now-smaller data range so the stored
# <autoFilter ref="..."/> shrinks along with the data instead of keeping the
# stale A1:G36 range.
$ws.AutoFilterMode = $false
$ws.Range("A1:G35").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the new
# filter range as well.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$35"
    }
}
